$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price column (D) cells whose new values could be
# misinterpreted as numbers by Excel, so the literal text (incl. trailing
# zeros / multi-dot formatted numbers) is preserved exactly as authored.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.224.85"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.425.37"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.42"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.01"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -1.78%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.74"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("E12").Value = "  +2.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.971.36"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000213"
$ws.Range("E14").Value = "  +4.76%  "

$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.43"
$ws.Range("E16").Value = "  -3.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.434.71"
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").Value = "  +3.92%  "

$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.282.11"
$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "464.37"
$ws.Range("E21").Value = "  +2.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.63"
$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("E23").Value = "  +3.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.41"
$ws.Range("E24").Value = "  +4.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.35"
$ws.Range("E25").Value = "  +19.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.30"
$ws.Range("E26").Value = "  +1.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "32.92"
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.71"
$ws.Range("E29").Value = "  +2.20%  "

$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("E31").Value = "  -3.37%  "

$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.90"
$ws.Range("E34").Value = "  -4.25%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.72"
$ws.Range("E36").Value = "  +8.46%  "

$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +4.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.326"
$ws.Range("E40").Value = "  +3.98%  "

$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("E42").Value = "  -0.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "144.66"
$ws.Range("E43").Value = "  +2.70%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  +10.10%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.07"
$ws.Range("E45").Value = "  +5.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.32"
$ws.Range("E46").Value = "  +2.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("E47").Value = "  +19.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.43"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.25"
$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0511"
$ws.Range("E50").Value = "  +25.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.12"
$ws.Range("E51").Value = "  +4.57%  "
